$d = $word.ActiveDocument

# Locate the "KEY ACHIEVEMENTS AND IMPACT" heading so we scope our edits to
# that section only (several bullet texts are duplicated earlier in the
# "PROFESSIONAL EXPERIENCE" section of this resume).
$sectionStart = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text.Trim() -eq "KEY ACHIEVEMENTS AND IMPACT") {
        $sectionStart = $i
        break
    }
}

if ($sectionStart -eq $null) {
    throw "Could not find KEY ACHIEVEMENTS AND IMPACT heading"
}

# Find the end of the section (the next Heading2-styled paragraph, i.e. the
# following top-level section, "TECHNICAL SKILLS").
$sectionEnd = $d.Paragraphs.Count
for ($i = $sectionStart + 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Style.NameLocal -eq "Heading 2") {
        $sectionEnd = $i - 1
        break
    }
}

# Within that section, rewrite the bullet paragraphs in place, one at a time,
# matching against their current text so we never touch the wrong paragraph.
for ($i = $sectionStart; $i -le $sectionEnd; $i++) {
    $para = $d.Paragraphs.Item($i)
    $text = $para.Range.Text

    if ($text -like "*Achieved 87% prediction accuracy for voter turnout*") {
        $para.Range.Text = "• Revenue generation: Delivered `$4.9M additional revenue through optimization"
    }
    elseif ($text -like "*Delivered `$4.9M additional revenue through continuous testing*") {
        $para.Range.Text = "• 23% conversion rate improvement"
    }
    elseif ($text -like "*Built redistricting platform used by thousands of analysts nationwide*") {
        $para.Range.Text = "• Executive authority: Briefed Presidents, Congressmen, Senators, Governors on election integrity, voter sentiment and postmortem analysis"
    }
    elseif ($text -like "*Trigonometric algorithm for boundary estimation*") {
        $para.Range.Text = "• Platform impact: Built redistricting system serving 12,847 analysts across 89 organizations"
    }
}

# Remove the two bullets that are dropped entirely, working from the bottom
# up within the section so earlier indices stay valid as paragraphs shift.
for ($i = $sectionEnd; $i -ge $sectionStart; $i--) {
    $para = $d.Paragraphs.Item($i)
    $text = $para.Range.Text

    if (($text -like "*Developed longitudinal data analysis methods using geospatial techniques*") -or
        ($text -like "*Discovered systematic race coding errors affecting all Black and Asian-American voters*")) {
        $para.Range.Delete()
    }
}
